{"js": "// Wrap a <w:body> inner-XML fragment into the full OOXML package payload\n// that Range.insertOoxml / insertOoxml expects.\nfunction wrapOoxml(bodyInnerXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + bodyInnerXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\n// Replace the whole paragraph that contains the first match of `needle`\n// with the raw `<w:p>\u2026</w:p>` fragment(s) supplied in `newParagraphsXml`.\nasync function replaceParagraphByText(context, needle, newParagraphsXml) {\n  const results = context.document.body.search(needle, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  const match = results.items[0];\n  const paragraph = match.paragraphs.getFirst();\n  const rng = paragraph.getRange(\"Whole\");\n  rng.insertOoxml(wrapOoxml(newParagraphsXml), Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst W = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\n// ---------------------------------------------------------------------------\n// 1. \"1. For Database\" heading becomes \"1. Environment \", followed by three\n//    new environment-info paragraphs, followed by the renumbered\n//    \"2. For Database\" heading.\n// ---------------------------------------------------------------------------\nconst block1 =\n  `<w:p ${W}>` +\n    '<w:pPr><w:pStyle w:val=\"Heading2\"/><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">1. Environment </w:t></w:r>' +\n  '</w:p>' +\n  `<w:p ${W}>` +\n    '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Java :</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> 17</w:t></w:r>' +\n  '</w:p>' +\n  `<w:p ${W}>` +\n    '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>SpringDoc</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>: 1.7.0</w:t></w:r>' +\n  '</w:p>' +\n  `<w:p ${W}>` +\n    '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Spring Boot 3.1.2</w:t></w:r>' +\n  '</w:p>' +\n  `<w:p ${W}>` +\n    '<w:pPr><w:pStyle w:val=\"Heading2\"/><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>2</w:t></w:r>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>. For Database</w:t></w:r>' +\n  '</w:p>';\n\nawait replaceParagraphByText(context, \"1. For Database\", block1);\n\n// ---------------------------------------------------------------------------\n// 2. \"2. Logging implementation\" -> split into \"3\" + \". Logging implementation\"\n// ---------------------------------------------------------------------------\nconst block2 =\n  `<w:p ${W}><w:pPr><w:pStyle w:val=\"Heading2\"/></w:pPr>` +\n  '<w:r><w:t>3</w:t></w:r>' +\n  '<w:r><w:t>. Logging implementation</w:t></w:r>' +\n  '</w:p>';\nawait replaceParagraphByText(context, \"2. Logging implementation\", block2);\n\n// ---------------------------------------------------------------------------\n// 3. \"3. Http response status \" -> split into \"4\" + \". Http response status \"\n// ---------------------------------------------------------------------------\nconst block3 =\n  `<w:p ${W}><w:pPr><w:pStyle w:val=\"Heading2\"/><w:tabs><w:tab w:val=\"right\" w:pos=\"9360\"/></w:tabs></w:pPr>` +\n  '<w:r><w:t>4</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">. Http response status </w:t></w:r>' +\n  '</w:p>';\nawait replaceParagraphByText(context, \"3. Http response status \", block3);\n\n// ---------------------------------------------------------------------------\n// 4. \"4. Unit test is planned to use Junit 5 and Mokito\" -> split into\n//    \"5\" + \". Unit test is planned to use Junit 5 and Mokito\"\n// ---------------------------------------------------------------------------\nconst block4 =\n  `<w:p ${W}><w:pPr><w:pStyle w:val=\"Heading2\"/><w:tabs><w:tab w:val=\"right\" w:pos=\"9360\"/></w:tabs></w:pPr>` +\n  '<w:r><w:t>5</w:t></w:r>' +\n  '<w:r><w:t>. Unit test is planned to use Junit 5 and Mokito</w:t></w:r>' +\n  '</w:p>';\nawait replaceParagraphByText(context, \"4. Unit test is planned to use Junit 5 and Mokito\", block4);\n\n// ---------------------------------------------------------------------------\n// 5. \"5. API E2E Browser Test\" -> becomes Heading 2, split into\n//    \"6\" + \". API E2E Browser Test\"\n// ---------------------------------------------------------------------------\nconst block5 =\n  `<w:p ${W}><w:pPr><w:pStyle w:val=\"Heading2\"/></w:pPr>` +\n  '<w:r><w:t>6</w:t></w:r>' +\n  '<w:r><w:t>. API E2E Browser Test</w:t></w:r>' +\n  '</w:p>';\nawait replaceParagraphByText(context, \"5. API E2E Browser Test\", block5);\n", "ps1": "$d = $word.ActiveDocument\n$wns = \"xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'\"\n\nfunction Get-ParagraphByText($needle) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($needle)\n    return $rng.Paragraphs(1)\n}\n\n# ---------------------------------------------------------------------------\n# 1. \"1. For Database\" heading becomes \"1. Environment \", followed by three\n#    new environment-info paragraphs, followed by the renumbered\n#    \"2. For Database\" heading.\n# ---------------------------------------------------------------------------\n$xml1 = @\"\n<w:p $wns>\n  <w:pPr><w:pStyle w:val='Heading2'/><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>\n  <w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>1. Environment </w:t></w:r>\n</w:p>\n<w:p $wns>\n  <w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>\n  <w:proofErr w:type='gramStart'/>\n  <w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Java :</w:t></w:r>\n  <w:proofErr w:type='gramEnd'/>\n  <w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> 17</w:t></w:r>\n</w:p>\n<w:p $wns>\n  <w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>\n  <w:proofErr w:type='spellStart'/>\n  <w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>SpringDoc</w:t></w:r>\n  <w:proofErr w:type='spellEnd'/>\n  <w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>: 1.7.0</w:t></w:r>\n</w:p>\n<w:p $wns>\n  <w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>\n  <w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Spring Boot 3.1.2</w:t></w:r>\n</w:p>\n<w:p $wns>\n  <w:pPr><w:pStyle w:val='Heading2'/><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>\n  <w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>2</w:t></w:r>\n  <w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>. For Database</w:t></w:r>\n</w:p>\n\"@\n\n$p1 = Get-ParagraphByText(\"1. For Database\")\n$p1.Range.InsertXML($xml1)\n\n# ---------------------------------------------------------------------------\n# 2. \"2. Logging implementation\" -> split into \"3\" + \". Logging implementation\"\n# ---------------------------------------------------------------------------\n$xml2 = \"<w:p $wns><w:pPr><w:pStyle w:val='Heading2'/></w:pPr><w:r><w:t>3</w:t></w:r><w:r><w:t>. Logging implementation</w:t></w:r></w:p>\"\n$p2 = Get-ParagraphByText(\"2. Logging implementation\")\n$p2.Range.InsertXML($xml2)\n\n# ---------------------------------------------------------------------------\n# 3. \"3. Http response status \" -> split into \"4\" + \". Http response status \"\n# ---------------------------------------------------------------------------\n$xml3 = \"<w:p $wns><w:pPr><w:pStyle w:val='Heading2'/><w:tabs><w:tab w:val='right' w:pos='9360'/></w:tabs></w:pPr><w:r><w:t>4</w:t></w:r><w:r><w:t xml:space='preserve'>. Http response status </w:t></w:r></w:p>\"\n$p3 = Get-ParagraphByText(\"3. Http response status \")\n$p3.Range.InsertXML($xml3)\n\n# ---------------------------------------------------------------------------\n# 4. \"4. Unit test is planned to use Junit 5 and Mokito\" -> split into\n#    \"5\" + \". Unit test is planned to use Junit 5 and Mokito\"\n# ---------------------------------------------------------------------------\n$xml4 = \"<w:p $wns><w:pPr><w:pStyle w:val='Heading2'/><w:tabs><w:tab w:val='right' w:pos='9360'/></w:tabs></w:pPr><w:r><w:t>5</w:t></w:r><w:r><w:t>. Unit test is planned to use Junit 5 and Mokito</w:t></w:r></w:p>\"\n$p4 = Get-ParagraphByText(\"4. Unit test is planned to use Junit 5 and Mokito\")\n$p4.Range.InsertXML($xml4)\n\n# ---------------------------------------------------------------------------\n# 5. \"5. API E2E Browser Test\" -> becomes Heading 2, split into\n#    \"6\" + \". API E2E Browser Test\"\n# ---------------------------------------------------------------------------\n$xml5 = \"<w:p $wns><w:pPr><w:pStyle w:val='Heading2'/></w:pPr><w:r><w:t>6</w:t></w:r><w:r><w:t>. API E2E Browser Test</w:t></w:r></w:p>\"\n$p5 = Get-ParagraphByText(\"5. API E2E Browser Test\")\n$p5.Range.InsertXML($xml5)\n"}
